$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 720, pushing existing rows 720..744 down to 721..745
$ws.Rows.Item(720).Insert()

# Populate the new row 720 with the new record
$ws.Cells.Item(720, 1).Value = 6
$ws.Cells.Item(720, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(720, 3).Value = "Metropolitana"
$ws.Cells.Item(720, 4).Value = 45075
$ws.Cells.Item(720, 5).Value = 13
$ws.Cells.Item(720, 6).Value = 100112044
$ws.Cells.Item(720, 7).Value = "Perejil"
$ws.Cells.Item(720, 8).Value = "Sin especificar"
$ws.Cells.Item(720, 9).Value = "Primera"
$ws.Cells.Item(720, 10).Value = 190
$ws.Cells.Item(720, 11).Value = 11000
$ws.Cells.Item(720, 12).Value = 12000
$ws.Cells.Item(720, 13).Value = 11421
$ws.Cells.Item(720, 14).Value = "$/docena de atados"
$ws.Cells.Item(720, 15).Value = "Región Metropolitana"
$ws.Cells.Item(720, 16).Value = 3807
$ws.Cells.Item(720, 17).Value = 3
$ws.Cells.Item(720, 18).Value = "Hortaliza"
